$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells for the team win/loss/tie record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from the last
# existing header cell onto the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in the team record for every data row (2 through 44)
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 71   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 91   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 1    # AF - Ties
}
